# Expansão das análises automáticas: adiciona colunas L, M, N
# (apoio_medio, contribuicoes, media_contribuicoes) ao resumo.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copia o estilo do cabeçalho existente (K1) para os novos cabeçalhos
# para que L1:N1 recebam o mesmo formato em negrito/borda.
$ws.Range("K1").Copy()
$ws.Range("L1:N1").PasteSpecial(-4122)

# Cabeçalhos das novas colunas
$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"

# Dados das novas colunas (linhas 2 a 7)
$ws.Range("L2").Value = 91.85312348260253
$ws.Range("M2").Value = 209535
$ws.Range("N2").Value = 322.3615384615385

$ws.Range("L3").Value = 89.17093558435907
$ws.Range("M3").Value = 54018
$ws.Range("N3").Value = 300.1

$ws.Range("L4").Value = 89.37434882498151
$ws.Range("M4").Value = 141221
$ws.Range("N4").Value = 132.6018779342723

$ws.Range("L5").Value = 91.95990423942952
$ws.Range("M5").Value = 62425
$ws.Range("N5").Value = 196.3050314465409

$ws.Range("L6").Value = 19.36290068160405
$ws.Range("M6").Value = 2129
$ws.Range("N6").Value = 15.65441176470588

$ws.Range("L7").Value = 24.85243295759227
$ws.Range("M7").Value = 79
$ws.Range("N7").Value = 4.9375
